$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C (Förändrad) for data rows 2 through 24:
# the date serial value changes from 46070 to 46072 for every row.
for ($row = 2; $row -le 24; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46070) {
        $cell.Value = 46072
    }
}
